$d = $word.ActiveDocument

# Find the paragraph containing "April 19" and replace the date text with [TODAYS DATE]
$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute("April 19, 2020", $true, $false, $false, $false, $false, $true, 1, $false, "[TODAYS DATE]", 2)
